# The document has a BTec logo picture (JPEG) in each header and a Pearson
# logo picture (PNG) in each footer. Both pairs of pictures currently carry
# mismatched sibling names ("image2.png"/"image1.jpg") - this renames them
# so header pictures are named "image2.jpg" and footer pictures are named
# "image1.png", matching the author's pairing fix.
$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            $shapeCount = $hdr.Range.InlineShapes.Count
            for ($j = 1; $j -le $shapeCount; $j++) {
                # Re-derive the InlineShapes collection from the shape's own
                # Range before writing - setting straight off the
                # HeaderFooter.Range collection can address a stale handle.
                $shp = $hdr.Range.InlineShapes.Item($j)
                $shp.Range.InlineShapes.Item(1).Name = "image2.jpg"
            }
        }

        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            $shapeCount = $ftr.Range.InlineShapes.Count
            for ($j = 1; $j -le $shapeCount; $j++) {
                $shp = $ftr.Range.InlineShapes.Item($j)
                $shp.Range.InlineShapes.Item(1).Name = "image1.png"
            }
        }
    }
}
